$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "version" column (C2:C119) from 0.8.5 to 0.10.0.
# This is a single shared string reused by every row, so the whole
# range is written in one shot to keep them pointing at one string.
$ws.Range("C2:C119").Value = "0.10.0"

# Update the per-row orch_session_issue_id GUIDs in column L.
$ws.Range("L2").Value = "7d3130df-50b7-429c-8354-be2781c2e78b"
$ws.Range("L3").Value = "970f3a23-308f-4875-aeed-a9c8010732b4"
$ws.Range("L4").Value = "2aa428f9-e8ca-4525-872f-bd6224121a31"
$ws.Range("L5").Value = "5f7614c3-9e73-47e7-9bc9-99999364ae19"
$ws.Range("L6").Value = "5f5046ce-d6bf-40f4-a113-90c87683310e"
$ws.Range("L7").Value = "0090848b-3389-4f3d-8899-9174dae0c622"
$ws.Range("L8").Value = "63c20d3e-5a88-4cd1-a6ba-491850a65151"
$ws.Range("L9").Value = "5c2f7272-f544-486c-8487-46b3ee126efc"
$ws.Range("L10").Value = "6516c754-1959-41ad-80c2-689a5a13341f"
$ws.Range("L11").Value = "f7bbe4fd-9832-4732-9014-dab8cc90433c"
$ws.Range("L12").Value = "79b706fc-a190-48de-a5be-38b9499c5ff4"
$ws.Range("L13").Value = "1bc36a0b-6bec-4d03-91b1-a77580353f2c"
$ws.Range("L14").Value = "73ee848c-b9fc-4ee4-af18-a6a5140092c9"
$ws.Range("L15").Value = "6d1a4149-7db8-479f-9b12-6a0a1ea06fa4"
$ws.Range("L16").Value = "8b14cf18-2ccd-4c7f-9146-d78f4d978131"
$ws.Range("L17").Value = "5485f038-d864-4f35-9576-c3fc3f94cc1c"
$ws.Range("L18").Value = "edcc23eb-1679-48d7-ab2f-32ddbf286025"
$ws.Range("L19").Value = "35cbf2f9-fb17-4ddc-881b-f249aad5c85a"
$ws.Range("L20").Value = "c03ca390-ee24-4c08-b0f2-8084f2ea28bb"
$ws.Range("L21").Value = "1524a424-a9a9-4dfa-89b1-00e74e3dc9d9"
$ws.Range("L22").Value = "e0f7d99a-d838-4aad-8f55-ca83037ca84d"
$ws.Range("L23").Value = "e3a8a36e-324a-4a21-8b3b-6b9e99e9bf47"
$ws.Range("L24").Value = "db33908d-84ab-48a0-8de9-a7344e2fdebc"
$ws.Range("L25").Value = "7f2bd43e-efad-4166-a197-caca396c7645"
$ws.Range("L26").Value = "edfad29c-cae4-45dc-a477-eb6d0b10520f"
$ws.Range("L27").Value = "52a046f3-f40d-4bed-a3c2-9d46ed164066"
$ws.Range("L28").Value = "1315191e-507c-4860-b5a8-49a930143651"
$ws.Range("L29").Value = "fd8de152-86a8-4f5a-91fc-0e8ad6c03a7e"
$ws.Range("L30").Value = "5a4e2085-7cdd-4aa2-8dc6-974c00ceb666"
$ws.Range("L31").Value = "6c00a0c2-4871-457a-a490-46698eea8d18"
$ws.Range("L32").Value = "408752d3-5563-40ee-a9d6-68f731bf48f9"
$ws.Range("L33").Value = "7b0932d3-3fcc-434f-a3a2-ecd786fa5d13"
$ws.Range("L34").Value = "b4a582bd-55b2-4988-b48a-c62c4b70a3d6"
$ws.Range("L35").Value = "155f15b6-89ba-41b1-bf39-b3e0cc337de7"
$ws.Range("L36").Value = "5e900cf0-fcd1-438e-8bcd-23dce6a0d248"
$ws.Range("L37").Value = "ef9164f5-5f5f-4a8c-a7cd-4ff1509a996b"
$ws.Range("L38").Value = "e93f8c06-7e5e-49b0-9e12-5c52fe1ebb85"
$ws.Range("L39").Value = "5bcdbb96-e63a-4f77-9bdf-767492bcfd44"
$ws.Range("L40").Value = "8bf8e88d-5c62-4921-a651-84a690d607f0"
$ws.Range("L41").Value = "025b8512-4c75-4d04-88d2-8f926b3ad26a"
$ws.Range("L42").Value = "76110869-eec6-43cc-a8a0-b8753a77d45c"
$ws.Range("L43").Value = "5b41f561-a1e0-4f99-8da6-80bbf546da39"
$ws.Range("L44").Value = "a1e06888-6e29-4858-b05c-64b3bfb14bbc"
$ws.Range("L45").Value = "d5d43fec-49d9-45de-b5a6-0d313887c5e5"
$ws.Range("L46").Value = "d91c129e-dc2c-45c6-812d-7919a5f38b7a"
$ws.Range("L47").Value = "067d8643-0b0f-4017-b27e-9b9223780315"
$ws.Range("L48").Value = "3f43598a-3824-42ec-ade3-f449cd4812a4"
$ws.Range("L49").Value = "e879406b-aa02-4284-ba23-b83f61df3720"
$ws.Range("L50").Value = "b9278cc8-f222-4f39-9200-65c99c03bd2d"
$ws.Range("L51").Value = "5c853da9-aba7-4495-9948-3f062fa3fac8"
$ws.Range("L52").Value = "2a90d219-3dcf-40cf-b85a-03ec2e035e32"
$ws.Range("L53").Value = "d52121f1-231f-4329-beab-1ca316bdcd13"
$ws.Range("L54").Value = "2b93b9ee-dd8e-475d-9183-03628907443b"
$ws.Range("L55").Value = "fcccaf8f-55bb-4430-b8bc-39e41c55ac01"
$ws.Range("L56").Value = "1eb8cd50-8dab-4b93-b30a-9c05a5b24696"
$ws.Range("L57").Value = "362a9ba3-b775-4b2c-a44b-ee0e0701fdb8"
$ws.Range("L58").Value = "e921b823-ce17-45d4-b389-1f7ea878370e"
$ws.Range("L59").Value = "500f86a2-08c5-42eb-9583-0d5812b7281f"
$ws.Range("L60").Value = "7103df56-e279-4012-a281-0e2114f9fa6f"
$ws.Range("L61").Value = "87be461b-dd68-499f-b9d5-adebef433033"
$ws.Range("L62").Value = "dd4081dd-3c95-4bb0-83bc-74afc385b6b2"
$ws.Range("L63").Value = "6a00c6da-9a2a-4aad-997f-f45be9997179"
$ws.Range("L64").Value = "4f142396-71c2-4a8d-ac64-3e0fc5043a09"
$ws.Range("L65").Value = "fdcba316-55b2-4385-9358-717507921100"
$ws.Range("L66").Value = "c6c7efbf-93c5-43f9-a778-b0de7a6a4d89"
$ws.Range("L67").Value = "3c3bea8b-a30e-411b-b2f6-84f738ba8a7a"
$ws.Range("L68").Value = "97bb2498-4268-4777-a8bc-10cc4974de88"
$ws.Range("L69").Value = "ad494fe4-0a9d-44bd-8813-524062768d08"
$ws.Range("L70").Value = "9fe506d8-6e9b-4ab6-9db6-449c814d229d"
$ws.Range("L71").Value = "b71c98ad-bc7c-49d5-87f7-4c3b6ccbb89c"
$ws.Range("L72").Value = "cf82cd2d-7973-4c20-ae2f-ae2655dbbc9a"
$ws.Range("L73").Value = "0fede892-aced-4d19-a33d-4ab95b3dc2aa"
$ws.Range("L74").Value = "0b328383-e38d-4efa-bfca-628edb277619"
$ws.Range("L75").Value = "b5298a09-0a07-411e-8ce8-2d94a760354d"
$ws.Range("L76").Value = "5ae8ea3d-e4dd-4d15-ad93-287aae58fe55"
$ws.Range("L77").Value = "a590e431-cb49-420c-bbf0-54269408bdff"
$ws.Range("L78").Value = "6e4de5b2-20e6-4a63-af23-c956d30a0db6"
$ws.Range("L79").Value = "b02edc1b-e15e-4b9e-b0b7-77f9ce798d64"
$ws.Range("L80").Value = "57c1bb4c-92fb-45fd-87e3-3d524dfd9ab0"
$ws.Range("L81").Value = "b78f0c90-163b-4ce2-95a6-f2e9d81ab56e"
$ws.Range("L82").Value = "dd5294bc-9ea8-4a84-95a2-ad95b403b325"
$ws.Range("L83").Value = "314f5962-e44d-4653-8973-6cdb4c1999c2"
$ws.Range("L84").Value = "1f8ee814-7404-4c62-b13d-b8d3c95bcceb"
$ws.Range("L85").Value = "5ca7f71c-47bd-414a-89bc-3328bee4b4f2"
$ws.Range("L86").Value = "8a01c7c5-0485-4e55-bb58-91396edbf552"
$ws.Range("L87").Value = "bd608a13-dc4e-44da-9deb-417739471840"
$ws.Range("L88").Value = "af37d7a2-20fd-458e-b4c4-8d1cbbb5efb3"
$ws.Range("L89").Value = "73b970dc-d0f2-4d06-8d7e-1c54f6c9a4f6"
$ws.Range("L90").Value = "cc5e3cbf-6172-4906-8b0b-3707b3b1f83d"
$ws.Range("L91").Value = "e2e3e84e-3239-4864-ac88-ac74797ed78a"
$ws.Range("L92").Value = "3863779d-d94e-4751-8149-4fc135bbdb5e"
$ws.Range("L93").Value = "aa47fcde-2351-4337-98b3-62c425a72a75"
$ws.Range("L94").Value = "9544db2e-42bc-4be7-a77e-03391b7aa215"
$ws.Range("L95").Value = "4a2eb67a-88fa-438f-8085-292f5c25efc7"
$ws.Range("L96").Value = "73034125-5056-4e70-964b-693962875c1d"
$ws.Range("L97").Value = "52644959-5383-466d-9694-20c3a3156de7"
$ws.Range("L98").Value = "0ada3647-1ff3-4aab-a783-9e5c19be20eb"
$ws.Range("L99").Value = "a1492ef7-b2bf-4d10-a71f-04217cd95c12"
$ws.Range("L100").Value = "eeb35705-edbe-474d-94b1-344809808c6c"
$ws.Range("L101").Value = "53bf2850-4371-4999-8cd1-0a8a18248f65"
$ws.Range("L102").Value = "6f4d082c-f89a-43e7-914c-165ce6a81828"
$ws.Range("L103").Value = "1842614e-2880-4dad-8a20-6dc676159e7b"
$ws.Range("L119").Value = "397cde6b-a7f0-4a29-9ab7-4519f9c41b80"
